$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting existing B:I to C:J
$ws.Columns("B:B").Insert()

# Set the header label for the new column
$ws.Range("B1").Value = "status_label"

# Map the colored-square status in column A to a French label in new column B
$statusMap = @{
    "vert"   = "🟩"
    "orange" = "🟧"
    "rouge"  = "🟥"
}

$rowCount = $ws.Range("A1").End(4).Row
for ($r = 2; $r -le $rowCount; $r++) {
    $statusValue = $ws.Cells.Item($r, 1).Value2
    foreach ($label in $statusMap.Keys) {
        if ($statusValue -eq $statusMap[$label]) {
            $ws.Cells.Item($r, 2).Value = $label
        }
    }
}
